$d = $word.ActiveDocument

# Target replacement for the first (and only) paragraph: split the last
# run so a gramStart/gramEnd proofErr pair wraps "no escopo ou seja", add
# a blank paragraph, then append two new paragraphs describing SCRUM /
# Product Owner / Scrum Master, with proofErr spell/gram markers matching
# the authored text.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t xml:space="preserve">Segundo a leitura do texto o modelo mais apropriado para este projeto é o modelo espiral pois como ainda </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:t>podem haver</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> mudanças no escopo do projeto e necessário ter a todo instante o </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">planejamento, a analise de riscos, a execução e a verificação. Caso não se encaixe </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:t>no escopo ou seja</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> alterado é possível refazer essas etapas.</w:t></w:r>' `
    + '</w:p>'

$paraBlank = '<w:p ' + $wNs + '/>'

$para2 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t xml:space="preserve">O SCRUM poderia ser aplicado pois o proprietário da empresa </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:t>esta</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> implementando esse projeto tendo em vista o aumento da produtividade, ou seja, os ROI.</w:t></w:r>' `
    + '</w:p>'

$para3 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t xml:space="preserve">Portanto o </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>propietario</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> seria o </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Product</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Owner</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> pois ele conhece a infraestrutura e as necessidades da empresa, o Scrum Master seria o responsável por gerir a equipe, ele é definido pela </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>propia</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> empresa e o time seria os programadores que irão desenvolver o sistema, responsáveis por definir metas e entregar o produto conforme os moldes solicitados</w:t></w:r>' `
    + '</w:p>'

$fullXml = $para1 + $paraBlank + $para2 + $para3

$target = $d.Paragraphs(1).Range
$target.InsertXML($fullXml)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
